$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1_2")

# 1. Edit cell D5 text
$ws.Range("D5").Value = "Vailand at Gum"

# 2. Add new row 8 with data
$ws.Range("A8").Value = 8291
$ws.Range("B8").Value = "Hum"
$ws.Range("C8").Value = "Tum"
$ws.Range("D8").Value = "LalaLand"

# 3. Set active cell selection to D6
$ws.Range("D6").Select()
